# Applies the "algorithm summary and roadmap updated" change:
# Inserts a new roadmap feature row ("Improve zoom function in image viewer")
# right above the "Add support for image derotation" row on the
# "Version 0.7.0" worksheet, pushing the remaining rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Version 0.7.0")

# The two screenshot pictures are anchored (by absolute position) next to the
# rows that will be pushed down by the insertion below. Remember their
# current vertical position so they can be nudged down by the height of the
# newly inserted row, keeping them aligned with the rows they illustrate.
$billede1 = $ws.Shapes.Item("Billede1")
$billede2 = $ws.Shapes.Item("Billede2")
$billede1Top = $billede1.Top
$billede2Top = $billede2.Top

# Insert a new row at row 27; this shifts rows 27..41 down to 28..42 and
# picks up the formatting of the row above it (row 26), which already uses
# the standard body style used throughout column A:G.
$ws.Rows("27:27").Insert()

# Fill in the new feature row.
$ws.Range("A27").Value = "Improve zoom function in image viewer"
$ws.Range("B27").Value = "Users criticize that there is no info shown about the current zoom status, and that there is no easy way to set the zoom status to 100%."
$ws.Range("C27").Value = "The zoom status should be included in one of the viewer corners (upper right?), e.g. ""125%"". Additionally, pressing ""1"" or double-clicking on the viewer image should reset to 100%."
$ws.Range("D27").Value = "Rolf"
$ws.Range("E27").Value = "Must have"
$ws.Range("F27").Value = "0.9.0"
$ws.Range("G27").Value = "open"

# Match the target custom row height used by similarly-sized entries.
$ws.Rows("27:27").RowHeight = 79.5

# Shift the two pictures down by the height of the newly inserted row so
# they keep pointing at the same feature rows as before.
$billede1.Top = $billede1Top + 79.5
$billede2.Top = $billede2Top + 79.5

# Update the selection to reflect where the editor last left the cursor.
$ws.Range("G27").Select()
